$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Raffa prova"
$ws.Range("B3").Value = "Stefano Pizzini | MediaserT"
$ws.Range("C3").Value = "Carlo Stedile | MAI UNA GIOIA"
$ws.Range("D3").Value = "MARCO HEIDEMPERGHER | U.S. Guarna"
$ws.Range("E3").Value = "Giovanni Torboli | F.C. Gorillaz"
$ws.Range("F3").Value = "Andrea Conzatti | FC Savignano"
